# Applies the "aggiornamento fino a 6 gennaio 2022" update: appends new
# daily COVID data rows (465-491) to the sheet, extending the used range
# from A1:D464 to A1:D491.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A (date) style/number-format down to the new rows by
# copying the format from the last existing data row (A464); values are
# overwritten individually below.
$ws.Range("A464").Copy($ws.Range("A465:A491"))

$data = @(
    ,@(465, 44539, 2, 14, 224.323025156225)
    ,@(466, 44540, 4, 14, 224.323025156225)
    ,@(467, 44541, 6, 20, 320.4614645088928)
    ,@(468, 44542, 0, 18, 288.4153180580035)
    ,@(469, 44543, 4, 20, 320.4614645088928)
    ,@(470, 44544, 0, 16, 256.3691716071143)
    ,@(471, 44545, 0, 16, 256.3691716071143)
    ,@(472, 44546, 0, 14, 224.323025156225)
    ,@(473, 44547, 3, 13, 208.2999519307803)
    ,@(474, 44548, 1, 8, 128.1845858035571)
    ,@(475, 44550, 5, 13, 208.2999519307803)
    ,@(476, 44551, 0, 9, 144.2076590290018)
    ,@(477, 44552, 0, 9, 144.2076590290018)
    ,@(478, 44553, 2, 11, 176.253805479891)
    ,@(479, 44554, 0, 11, 176.253805479891)
    ,@(480, 44555, 3, 11, 176.253805479891)
    ,@(481, 44556, 6, 16, 256.3691716071143)
    ,@(482, 44557, 6, 17, 272.3922448325589)
    ,@(483, 44558, 11, 28, 448.6460503124499)
    ,@(484, 44559, 2, 30, 480.6921967633392)
    ,@(485, 44560, 2, 30, 480.6921967633392)
    ,@(486, 44561, 6, 36, 576.8306361160071)
    ,@(487, 44562, 4, 37, 592.8537093414517)
    ,@(488, 44563, 17, 48, 769.1075148213428)
    ,@(489, 44564, 17, 59, 945.3613203012337)
    ,@(490, 44565, 3, 51, 817.1767344976765)
    ,@(491, 44566, 17, 66, 1057.522832879346)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}
